$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player/position/team data (rows 2-19), reflecting the reordered roster
$data = @(
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Dillon Brooks", "SG,SF", "Houston Rockets"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Draymond Green", "PF,C", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
